# Updates the crypto price/volume snapshot on Sheet1 (cryptos.xlsx).
# Cell values that look like plain decimal numbers (e.g. "63.63") are
# written with a leading apostrophe so Excel stores them as literal text
# (matching the workbook's existing text-formatted Price column) instead
# of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.442.88"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.163.95"
$ws.Range("E3").Value = "  +2.89%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'228.62"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").Value = "'63.63"
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").Value = "'0.0854"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'16.08"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "2.485.46"
$ws.Range("E13").Value = "  +2.99%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'0.814"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "2.170.67"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "39.473.61"
$ws.Range("D19").Value = "'6.23"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").Value = "'71.86"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "'229.54"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.34"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'172.05"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.52"
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("D29").Value = "'19.86"
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("D30").Value = "'1.42"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "'2.69"
$ws.Range("E31").Value = "  +6.51%  "
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'4.63"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").Value = "'4.72"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "'7.06"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "'3.61"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").Value = "'103.11"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").Value = "'17.81"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("D43").Value = "1.521.34"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  +3.92%  "
$ws.Range("D45").Value = "'1.11"
$ws.Range("E45").Value = "  +5.60%  "
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").Value = "'0.0926"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").Value = "'4.26"
$ws.Range("E48").Value = "  +3.21%  "
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "2.369.09"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("E51").Value = "  -0.73%  "
